$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.755.96"
$ws.Range("D3").Value = "1.849.78"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4269"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07325"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.821.69"
$ws.Range("E13").Value = "  -3.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.583"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.355"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06927"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("E18").Value = "  -2.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008912"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.61%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").Value = "27.757.80"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.993"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("D25").Value = "2.088.18"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.966"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.252"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.906"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08947"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7706"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.969"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.106"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.05%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05401"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.84%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.097"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01959"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.77%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5125"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.915"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.291"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06603"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4772"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.637"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.40%  "
